# Regenerate the localization-status report:
#   - the "Status" for the zh-cn / de-de handoff moved on from
#     "Ready for handoff" to "In Translation"
#   - the Status columns are narrower now that the new text is shorter,
#     so re-fit their widths.

$wb = $excel.ActiveWorkbook

# --- 1. Update the Status text everywhere it appears -----------------
# Overview sheet: columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F2").Value = "In Translation"

# Per-locale detail sheets: column C (Status), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Re-fit the Status columns to the new (shorter) text -----------
# ColumnWidth = 12.5 lands the stored column width right on the same
# pixel grid cell that Excel's own AutoFit would use for this text.
$wsOverview.Columns.Item("E").ColumnWidth = 12.5
$wsOverview.Columns.Item("F").ColumnWidth = 12.5

$wsZhCn.Columns.Item("C").ColumnWidth = 12.5
$wsDeDe.Columns.Item("C").ColumnWidth = 12.5
